$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so the exact string (not a converted number)
# is stored, matching the original inline-string cell contents.
$cells = @{
    "D2" = "297.02"
    "E2" = "1.01%"
    "D3" = "31.61"
    "E3" = "1.89%"
    "D4" = "4.973"
    "E4" = "0.44%"
    "D5" = "0.07678"
    "E5" = "4.59%"
    "D6" = "2.250"
    "E6" = "-1.87%"
    "D7" = "7.869"
    "E7" = "1.85%"
    "D8" = "0.9244"
    "E8" = "1.74%"
    "D9" = "0.09796"
    "E9" = "22.28%"
    "D10" = "0.1746"
    "E10" = "4.01%"
    "D11" = "0.08402"
    "E11" = "3.82%"
    "D12" = "0.03249"
    "E12" = "4.68%"
    "D13" = "0.09827"
    "E13" = "-2.50%"
    "D14" = "0.001480"
    "E14" = "-2.63%"
    "D15" = "0.005763"
    "E15" = "-0.68%"
    "D16" = "3.522"
    "E16" = "0.90%"
    "D17" = "3.788"
    "E17" = "1.22%"
    "D18" = "2.196"
    "E18" = "5.78%"
    "D19" = "0.3360"
    "E19" = "0.88%"
    "D20" = "0.1315"
    "E20" = "0.87%"
    "D21" = "4.073"
    "E21" = "2.21%"
    "D22" = "0.2277"
    "E22" = "8.38%"
    "D23" = "0.04499"
    "E23" = "-1.16%"
    "D24" = "0.001213"
    "E24" = "-0.06%"
    "D25" = "0.004360"
    "E25" = "-6.29%"
    "D26" = "0.0001290"
    "E26" = "-0.85%"
    "D27" = "0.0003367"
    "E27" = "-0.86%"
    "D39" = "0.01698"
    "E39" = "5.75%"
    "D40" = "0.04631"
    "E40" = "4.64%"
    "D41" = "0.007537"
    "D42" = "0.009738"
    "E42" = "12.82%"
    "D43" = "0.1386"
    "E43" = "4.24%"
    "D44" = "0.002055"
    "E44" = "0.18%"
    "D45" = "0.01046"
    "E45" = "9.87%"
    "D46" = "0.00006026"
    "E46" = "1.73%"
    "E47" = "-0.86%"
    "D48" = "2.551"
    "E48" = "13.83%"
    "E49" = "-31.59%"
    "D50" = "0.00002083"
    "E50" = "-0.86%"
    "D51" = "0.0001984"
    "E51" = "-0.86%"
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}
